$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A8").Value = 3
$ws.Range("B8").Value = "Стандарт"
$ws.Range("C8").Value = 790
$ws.Range("D8").Value = 890
$ws.Range("E8").Value = 140

$ws.Range("A10").Select()
